$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New metric rows to append starting at row 107.
# Each metric has 5 tier rows: Lowest, Low, Medium, High, Highest.
$data = @(
    @("temperatureAvgJan", "Lowest", 0, 0),
    @("temperatureAvgJan", "Low", 0, 0),
    @("temperatureAvgJan", "Medium", 0, 0),
    @("temperatureAvgJan", "High", 0, 0),
    @("temperatureAvgJan", "Highest", 0, 0),
    @("waveMax*CMin*C", "Lowest", 0, 0.5),
    @("waveMax*CMin*C", "Low", 0.5, 5),
    @("waveMax*CMin*C", "Medium", 5, 10),
    @("waveMax*CMin*C", "High", 10, 40),
    @("waveMax*CMin*C", "Highest", 40, 9999),
    @("waveMax95Pct", "Lowest", 0, 1),
    @("waveMax95Pct", "Low", 1, 5),
    @("waveMax95Pct", "Medium", 5, 10),
    @("waveMax95Pct", "High", 10, 15),
    @("waveMax95Pct", "Highest", 15, 9999)
)

$startRow = 107
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
